$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.093498229980469
$ws.Range("B1").Value = 3.066576957702637
$ws.Range("C1").Value = 2.392396450042725
$ws.Range("D1").Value = 2.223505735397339
$ws.Range("E1").Value = 1.880835056304932
